$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Range, [string]$Value)
    $Range.NumberFormat = "@"
    $Range.Value = $Value
    $Range.Style = "Normal"
}

$ws.Range("D2").Value = "68.992.12"
$ws.Range("E2").Value = "  +0.35%  "
$ws.Range("D3").Value = "3.747.90"
$ws.Range("E3").Value = "  -0.59%  "
Set-TextValue $ws.Range("D4") "0.999"
$ws.Range("E4").Value = "  -0.12%  "
Set-TextValue $ws.Range("D5") "602.36"
$ws.Range("E5").Value = "  -0.12%  "
Set-TextValue $ws.Range("D6") "165.52"
$ws.Range("E6").Value = "  -2.42%  "
$ws.Range("D7").Value = "3.746.95"
$ws.Range("E7").Value = "  -0.53%  "
$ws.Range("E8").Value = "  -0.04%  "
$ws.Range("E9").Value = "  +0.29%  "
Set-TextValue $ws.Range("D10") "0.172"
$ws.Range("E10").Value = "  +4.72%  "
Set-TextValue $ws.Range("D11") "6.38"
$ws.Range("E11").Value = "  +0.53%  "
$ws.Range("E12").Value = "  -1.07%  "
Set-TextValue $ws.Range("D13") "37.67"
$ws.Range("E13").Value = "  -2.04%  "
Set-TextValue $ws.Range("D14") "0.0000248"
$ws.Range("E14").Value = "  +0.43%  "
$ws.Range("D15").Value = "4.374.40"
$ws.Range("E15").Value = "  -0.59%  "
$ws.Range("D16").Value = "3.758.04"
$ws.Range("E16").Value = "  -0.14%  "
$ws.Range("D17").Value = "68.934.75"
$ws.Range("E17").Value = "  +0.18%  "
$ws.Range("E18").Value = "  +1.20%  "
Set-TextValue $ws.Range("D19") "17.74"
$ws.Range("E19").Value = "  +3.01%  "
$ws.Range("E20").Value = "  -0.95%  "
Set-TextValue $ws.Range("D21") "11.21"
$ws.Range("E21").Value = "  +3.64%  "
Set-TextValue $ws.Range("D22") "489.99"
$ws.Range("E22").Value = "  -1.28%  "
$ws.Range("E23").Value = "  -1.25%  "
$ws.Range("B24").Value = "PEPE"
$ws.Range("C24").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
Set-TextValue $ws.Range("D24") "0.0000149"
$ws.Range("E24").Value = "  +2.51%  "
$ws.Range("B25").Value = "Litecoin"
$ws.Range("C25").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
Set-TextValue $ws.Range("D25") "84.70"
$ws.Range("E25").Value = "  -1.24%  "
$ws.Range("E26").Value = "  -2.46%  "
$ws.Range("E27").Value = "  -1.06%  "
Set-TextValue $ws.Range("D28") "10.04"
$ws.Range("E28").Value = "  -2.08%  "
$ws.Range("E29").Value = "  -0.08%  "
$ws.Range("E30").Value = "  -0.61%  "
$ws.Range("E31").Value = "  +1.76%  "
Set-TextValue $ws.Range("D32") "2.43"
$ws.Range("E32").Value = "  -4.40%  "
$ws.Range("D33").Value = "3.895.78"
$ws.Range("E33").Value = "  -0.51%  "
Set-TextValue $ws.Range("D34") "31.57"
$ws.Range("E34").Value = "  -2.00%  "
$ws.Range("D35").Value = "3.682.90"
$ws.Range("E35").Value = "  -0.55%  "
Set-TextValue $ws.Range("D36") "0.108"
$ws.Range("E36").Value = "  -1.11%  "
Set-TextValue $ws.Range("D37") "5.93"
$ws.Range("E37").Value = "  +1.04%  "
$ws.Range("E38").Value = "  -1.12%  "
$ws.Range("E39").Value = "  +3.98%  "
$ws.Range("E41").Value = "  +7.75%  "
$ws.Range("E42").Value = "  -0.46%  "
Set-TextValue $ws.Range("D43") "48.56"
$ws.Range("E43").Value = "  -0.74%  "
$ws.Range("E44").Value = "  +0.20%  "
Set-TextValue $ws.Range("D45") "423.72"
$ws.Range("E45").Value = "  -3.96%  "
$ws.Range("E46").Value = "  -1.16%  "
$ws.Range("E47").Value = "  +0.00%  "
Set-TextValue $ws.Range("D48") "40.08"
$ws.Range("E48").Value = "  -1.88%  "
Set-TextValue $ws.Range("D49") "141.76"
$ws.Range("E49").Value = "  +0.44%  "
Set-TextValue $ws.Range("D50") "1.31"
$ws.Range("E50").Value = "  +7.70%  "
$ws.Range("D51").Value = "2.780.40"
$ws.Range("E51").Value = "  -1.94%  "
